$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header for column C
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fund rows data (name, old B value, new C value) in the new order (rows 2..34)
$funds = @(
    @("1810 Renta variable", 605392.62, 640930.23),
    @("1822 Raices Valores Negociables", 990529.11, 988874.8199999999),
    @("Adcap Balanceado V", 295014.52, 266099.34),
    @("Adcap IOL Acciones Argentina", 19732, 19635.72),
    @("Adcap Wise", 532629.02, 477053.67),
    @("Allaria Acciones", 78538.86, 78354.11),
    @("Alpha Acciones", 389064.75, 390001.9),
    @("Alpha Mega", 1389912.42, 1389806.51),
    @("Alpha renta balan global", 1053603.04, 1048935.77),
    @("Argenfunds", 38732.12, 38634.65),
    @("Arpenta acciones", 5529.2, 5520.34),
    @("Balanz", 8645.709999999999, 8929.639999999999),
    @("Delta gestion V", 2578167.38, 2322965.99),
    @("FBA Acciones Argentinas", 166921.24, 167527.66),
    @("FBA Calificado", 160622.95, 160514.35),
    @("Fima Acciones", 1624102.04, 1745608.11),
    @("Fima PB Acciones", 201704.83, 200380.02),
    @("Goal Acciones Argentinas", 81289.53999999999, 81638.42),
    @("Goal acciones plus", 10466.79, 10519.39),
    @("HF Acciones Argentinas", 225352.3, 226172.57),
    @("HF Acciones Lideres", 373067.04, 373678.79),
    @("IAM Renta Variable", 124895.09, 131053.92),
    @("IEB Value", 15370.68, 15379.56),
    @("Lombardi", 124539.94, 138404.74),
    @("MAF", 740.3099999999999, 731.4),
    @("Megainver", 103135.9, 103046.6),
    @("Pellegrini Acciones", 303090.71, 302296.96),
    @("Pionero Acciones", 350780.68, 349269),
    @("Premier Renta Variable", 16033.92, 15636.3),
    @("Quinquela Acciones", 327476.6, 328656.13),
    @("Rofex 20 Renta Variable", 245846.51, 245588.31),
    @("Supefondo RV", 3197643.34, 3616847.15),
    @("Toronto Trust Multimercado", 100109.46, 100136.78)
)

$row = 2
foreach ($fund in $funds) {
    $ws.Cells.Item($row, 1).Value = $fund[0]
    $ws.Cells.Item($row, 2).Value = $fund[1]
    $ws.Cells.Item($row, 3).Value = $fund[2]
    $row++
}

# avg row (35)
$ws.Cells.Item(35, 1).Value = "avg"
$ws.Cells.Item(35, 2).Value = 476929.72
$ws.Cells.Item(35, 3).Value = 484509.97

# total row (36)
$ws.Cells.Item(36, 1).Value = "total"
$ws.Cells.Item(36, 2).Value = 15738680.62
$ws.Cells.Item(36, 3).Value = 15988828.85
